# Update countries & provincias Spain
# Refresh COVID stats (timestamp 16:28 -> 17:45) and re-sort-driven row shifts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 17:45"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7953173
$ws.Range("C4").Value = 7668
$ws.Range("D4").Value = 5090782
$ws.Range("E4").Value = 2643076
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 219315

# Row 5: India
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = 7063955
$ws.Range("C5").Value = 12412
$ws.Range("D5").Value = 6087588
$ws.Range("E5").Value = 867955
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 41
$ws.Range("H5").Value = 108412

# Row 15: Reino Unido
$ws.Range("A15").Value = "Reino Unido"
$ws.Range("B15").Value = 603716
$ws.Range("C15").Value = 12872
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 65
$ws.Range("H15").Value = 42825

# Row 20: Italia
$ws.Range("A20").Value = "Italia"
$ws.Range("B20").Value = 354950
$ws.Range("C20").Value = 5456
$ws.Range("D20").Value = 239709
$ws.Range("E20").Value = 79075
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 36166

# Row 25: Alemania
$ws.Range("A25").Value = "Alemania"
$ws.Range("B25").Value = 324596
$ws.Range("C25").Value = 1143
$ws.Range("D25").Value = 273500
$ws.Range("E25").Value = 41399
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 9697

# Row 39: Republica Dominicana
$ws.Range("A39").Value = "Republica Dominicana"
$ws.Range("B39").Value = 118477
$ws.Range("C39").Value = 463
$ws.Range("D39").Value = 94084
$ws.Range("E39").Value = 22220
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 2173

# Row 48: Guatemala
$ws.Range("A48").Value = "Guatemala"
$ws.Range("B48").Value = 97715
$ws.Range("C48").Value = 171
$ws.Range("D48").Value = 86582
$ws.Range("E48").Value = 7749
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 19
$ws.Range("H48").Value = 3384

# Row 49: Japon
$ws.Range("A49").Value = "Japon"
$ws.Range("B49").Value = 88912
$ws.Range("C49").Value = 679
$ws.Range("D49").Value = 81824
$ws.Range("E49").Value = 5461
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1627

# Row 62: Singapur
$ws.Range("A62").Value = "Singapur"
$ws.Range("B62").Value = 57876
$ws.Range("C62").Value = 10
$ws.Range("D62").Value = 57705
$ws.Range("E62").Value = 144
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 27

# Row 84: Jordania
$ws.Range("A84").Value = "Jordania"
$ws.Range("B84").Value = 24926
$ws.Range("C84").Value = 928
$ws.Range("D84").Value = 6101
$ws.Range("E84").Value = 18634
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 10
$ws.Range("H84").Value = 191

# Row 85: Corea del Sur
$ws.Range("A85").Value = "Corea del Sur"
$ws.Range("B85").Value = 24606
$ws.Range("C85").Value = 58
$ws.Range("D85").Value = 22693
$ws.Range("E85").Value = 1481
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 432

# Row 86: Bulgaria
$ws.Range("A86").Value = "Bulgaria"
$ws.Range("B86").Value = 24319
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 15818
$ws.Range("E86").Value = 7610
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 891

# Row 97: Albania
$ws.Range("A97").Value = "Albania"
$ws.Range("B97").Value = 15399
$ws.Range("C97").Value = 168
$ws.Range("D97").Value = 9500
$ws.Range("E97").Value = 5479
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 4
$ws.Range("H97").Value = 420

# Row 98: Senegal
$ws.Range("A98").Value = "Senegal"
$ws.Range("B98").Value = 15268
$ws.Range("C98").Value = 24
$ws.Range("D98").Value = 13297
$ws.Range("E98").Value = 1657
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 314

# Row 99: Montenegro
$ws.Range("A99").Value = "Montenegro"
$ws.Range("B99").Value = 13869
$ws.Range("C99").Value = 228
$ws.Range("D99").Value = 9620
$ws.Range("E99").Value = 4047
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 4
$ws.Range("H99").Value = 202

# Row 100: Sudan
$ws.Range("A100").Value = "Sudan"
$ws.Range("B100").Value = 13670
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 6764
$ws.Range("E100").Value = 6070
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 836

# Row 104: Guinea
$ws.Range("A104").Value = "Guinea"
$ws.Range("B104").Value = 11022
$ws.Range("C104").Value = 26
$ws.Range("D104").Value = 10324
$ws.Range("E104").Value = 629
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 69

# Row 105: Maldivas
$ws.Range("A105").Value = "Maldivas"
$ws.Range("B105").Value = 10859
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 9683
$ws.Range("E105").Value = 1141
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 35

# Row 110: Uganda
$ws.Range("A110").Value = "Uganda"
$ws.Range("B110").Value = 9801
$ws.Range("C110").Value = 100
$ws.Range("D110").Value = 6109
$ws.Range("E110").Value = 3599
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 93

# Row 111: Luxemburgo
$ws.Range("A111").Value = "Luxemburgo"
$ws.Range("B111").Value = 9722
$ws.Range("C111").Value = 208
$ws.Range("D111").Value = 8038
$ws.Range("E111").Value = 1552
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 132

# Row 116: Jamaica
$ws.Range("A116").Value = "Jamaica"
$ws.Range("B116").Value = 7718
$ws.Range("C116").Value = 159
$ws.Range("D116").Value = 3162
$ws.Range("E116").Value = 4417
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 139

# Row 122: Cuba
$ws.Range("A122").Value = "Cuba"
$ws.Range("B122").Value = 5978
$ws.Range("C122").Value = 30
$ws.Range("D122").Value = 5540
$ws.Range("E122").Value = 315
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 123

# Row 123: Malaui
$ws.Range("A123").Value = "Malaui"
$ws.Range("B123").Value = 5821
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 4647
$ws.Range("E123").Value = 994
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 180

# Row 135: Sri Lanka
$ws.Range("A135").Value = "Sri Lanka"
$ws.Range("B135").Value = 4689
$ws.Range("C135").Value = 61
$ws.Range("D135").Value = 3307
$ws.Range("E135").Value = 1369
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 13

# Row 136: Siria
$ws.Range("A136").Value = "Siria"
$ws.Range("B136").Value = 4673
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 1271
$ws.Range("E136").Value = 3181
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 221

# Row 144: Gambia
$ws.Range("A144").Value = "Gambia"
$ws.Range("B144").Value = 3632
$ws.Range("C144").Value = 4
$ws.Range("D144").Value = 2543
$ws.Range("E144").Value = 972
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 117

# Row 149: Sudan del Sur
$ws.Range("A149").Value = "Sudan del Sur"
$ws.Range("B149").Value = 2777
$ws.Range("C149").Value = 16
$ws.Range("D149").Value = 1290
$ws.Range("E149").Value = 1432
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 55

# Row 164: Lesoto
$ws.Range("A164").Value = "Lesoto"
$ws.Range("B164").Value = 1805
$ws.Range("C164").Value = 5
$ws.Range("D164").Value = 961
$ws.Range("E164").Value = 802
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 42
